## Generate Report for Handoff
## The "b.md" localization row has been handed off again: its status moves
## from "Handed back: in sync with en-US" to "Ready for handoff", the
## handoff file names are rev'd to the new content hash, and the handoff
## timestamps are updated.

$wb = $excel.ActiveWorkbook

$Overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# 1. Overview sheet - row 3 is the "b.md" file.
# ---------------------------------------------------------------------
$Overview.Range("B3").Value = "Ready for handoff"
$Overview.Range("C3").Value = "Ready for handoff"
$Overview.Range("D3").Value = "2016-22-14 02:22:44"

# ---------------------------------------------------------------------
# 2. zh-cn sheet - row 3 is the "b.md" file.
# ---------------------------------------------------------------------
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-14 02:22:42"

foreach ($link in $zhcn.Hyperlinks) {
    if ($link.Range.Address() -eq '$D$3') {
        $link.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# 3. de-de sheet - row 3 is the "b.md" file.
# ---------------------------------------------------------------------
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-14 02:22:44"

foreach ($link in $dede.Hyperlinks) {
    if ($link.Range.Address() -eq '$D$3') {
        $link.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
